$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "26-22="
$t.Cell(1, 2).Range.Text = "87-51="
$t.Cell(1, 3).Range.Text = "23+40="
$t.Cell(1, 4).Range.Text = "57+28="
$t.Cell(1, 5).Range.Text = "82+9="
$t.Cell(2, 1).Range.Text = "25+31="
$t.Cell(2, 2).Range.Text = "8+78="
$t.Cell(2, 3).Range.Text = "72-26="
$t.Cell(2, 4).Range.Text = "86-34="
$t.Cell(2, 5).Range.Text = "52-38="
$t.Cell(3, 1).Range.Text = "96-14="
$t.Cell(3, 2).Range.Text = "28+68="
$t.Cell(3, 3).Range.Text = "8+38="
$t.Cell(3, 4).Range.Text = "95-90="
$t.Cell(3, 5).Range.Text = "23-14="
$t.Cell(4, 1).Range.Text = "55+12="
$t.Cell(4, 2).Range.Text = "7+61="
$t.Cell(4, 3).Range.Text = "55+15="
$t.Cell(4, 4).Range.Text = "28+51="
$t.Cell(4, 5).Range.Text = "38-20="
$t.Cell(5, 1).Range.Text = "64+3="
$t.Cell(5, 2).Range.Text = "14+62="
$t.Cell(5, 3).Range.Text = "46-20="
$t.Cell(5, 4).Range.Text = "42-24="
$t.Cell(5, 5).Range.Text = "70-13="
$t.Cell(6, 1).Range.Text = "56+12="
$t.Cell(6, 2).Range.Text = "23-3="
$t.Cell(6, 3).Range.Text = "39-25="
$t.Cell(6, 4).Range.Text = "37+26="
$t.Cell(6, 5).Range.Text = "84-59="
$t.Cell(7, 1).Range.Text = "74-56="
$t.Cell(7, 2).Range.Text = "85-32="
$t.Cell(7, 3).Range.Text = "29-5="
$t.Cell(7, 4).Range.Text = "32-11="
$t.Cell(7, 5).Range.Text = "48-44="
$t.Cell(8, 1).Range.Text = "84-81="
$t.Cell(8, 2).Range.Text = "36+39="
$t.Cell(8, 3).Range.Text = "55+3="
$t.Cell(8, 4).Range.Text = "71-31="
$t.Cell(8, 5).Range.Text = "40+50="
$t.Cell(9, 1).Range.Text = "60+22="
$t.Cell(9, 2).Range.Text = "7+59="
$t.Cell(9, 3).Range.Text = "30+23="
$t.Cell(9, 4).Range.Text = "10+33="
$t.Cell(9, 5).Range.Text = "24+9="
$t.Cell(10, 1).Range.Text = "30+17="
$t.Cell(10, 2).Range.Text = "12+7="
$t.Cell(10, 3).Range.Text = "17+39="
$t.Cell(10, 4).Range.Text = "84-23="
$t.Cell(10, 5).Range.Text = "35-34="
$t.Cell(11, 1).Range.Text = "80-32="
$t.Cell(11, 2).Range.Text = "28-12="
$t.Cell(11, 3).Range.Text = "58+32="
$t.Cell(11, 4).Range.Text = "39-10="
$t.Cell(11, 5).Range.Text = "74+12="
$t.Cell(12, 1).Range.Text = "88-1="
$t.Cell(12, 2).Range.Text = "82-34="
$t.Cell(12, 3).Range.Text = "0+10="
$t.Cell(12, 4).Range.Text = "0+5="
$t.Cell(12, 5).Range.Text = "68-34="
$t.Cell(13, 1).Range.Text = "49-38="
$t.Cell(13, 2).Range.Text = "7+18="
$t.Cell(13, 3).Range.Text = "6+68="
$t.Cell(13, 4).Range.Text = "82-9="
$t.Cell(13, 5).Range.Text = "10-7="
$t.Cell(14, 1).Range.Text = "38-13="
$t.Cell(14, 2).Range.Text = "93-47="
$t.Cell(14, 3).Range.Text = "53+23="
$t.Cell(14, 4).Range.Text = "2+50="
$t.Cell(14, 5).Range.Text = "23+34="
$t.Cell(15, 1).Range.Text = "35+62="
$t.Cell(15, 2).Range.Text = "93+6="
$t.Cell(15, 3).Range.Text = "18+16="
$t.Cell(15, 4).Range.Text = "40-11="
$t.Cell(15, 5).Range.Text = "73-4="
$t.Cell(16, 1).Range.Text = "68-63="
$t.Cell(16, 2).Range.Text = "15+74="
$t.Cell(16, 3).Range.Text = "2+65="
$t.Cell(16, 4).Range.Text = "81-71="
$t.Cell(16, 5).Range.Text = "11+45="
$t.Cell(17, 1).Range.Text = "93-37="
$t.Cell(17, 2).Range.Text = "30-14="
$t.Cell(17, 3).Range.Text = "82-39="
$t.Cell(17, 4).Range.Text = "99-67="
$t.Cell(17, 5).Range.Text = "24-23="
$t.Cell(18, 1).Range.Text = "16+30="
$t.Cell(18, 2).Range.Text = "48+39="
$t.Cell(18, 3).Range.Text = "18+25="
$t.Cell(18, 4).Range.Text = "59-25="
$t.Cell(18, 5).Range.Text = "44-25="
$t.Cell(19, 1).Range.Text = "87-49="
$t.Cell(19, 2).Range.Text = "79-31="
$t.Cell(19, 3).Range.Text = "71-59="
$t.Cell(19, 4).Range.Text = "36-11="
$t.Cell(19, 5).Range.Text = "98-0="
$t.Cell(20, 1).Range.Text = "40-27="
$t.Cell(20, 2).Range.Text = "40-21="
$t.Cell(20, 3).Range.Text = "50-10="
$t.Cell(20, 4).Range.Text = "28-25="
$t.Cell(20, 5).Range.Text = "40+58="
